$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NATMI LR-pair computed columns (M:T) with refreshed TPM-derived values.
# N = ReceptorExpressingCells * AvgExpr; O/P = specificity of avg/total receptor expr across target clusters;
# Q/R = ligand*receptor avg/total expr weight; S/T = ligand*receptor derived specificity product.

$ws.Range("M2").Value = 3.867218333333334
$ws.Range("N2").Value = 11.601655
$ws.Range("O2").Value = 0.1566152977872902
$ws.Range("P2").Value = 0.1566152977872902
$ws.Range("Q2").Value = 0.6183695005727778
$ws.Range("R2").Value = 5.565325505155
$ws.Range("S2").Value = 0.004146815913990309
$ws.Range("T2").Value = 0.004146815913990309

$ws.Range("N3").Value = 33.813685
$ws.Range("O3").Value = 0.4564642152831324
$ws.Range("P3").Value = 0.4564642152831324
$ws.Range("S3").Value = 0.01208613142423692
$ws.Range("T3").Value = 0.01208613142423692

$ws.Range("M4").Value = 5.654344666666667
$ws.Range("N4").Value = 16.963034
$ws.Range("O4").Value = 0.2289906587711778
$ws.Range("P4").Value = 0.2289906587711778
$ws.Range("Q4").Value = 0.9041315969815557
$ws.Range("R4").Value = 8.137184372834
$ws.Range("S4").Value = 0.006063150416105174
$ws.Range("T4").Value = 0.006063150416105175

$ws.Range("M5").Value = 0.819389
$ws.Range("N5").Value = 2.458167
$ws.Range("O5").Value = 0.03318376186120772
$ws.Range("P5").Value = 0.03318376186120772
$ws.Range("Q5").Value = 0.1310205742296667
$ws.Range("R5").Value = 1.179185168067
$ws.Range("S5").Value = 0.0008786303363482032
$ws.Range("T5").Value = 0.0008786303363482033

$ws.Range("M6").Value = 3.080288333333333
$ws.Range("N6").Value = 9.240865
$ws.Range("O6").Value = 0.1247460662971919
$ws.Range("P6").Value = 0.1247460662971919
$ws.Range("Q6").Value = 0.4925391312627778
$ws.Range("R6").Value = 4.432852181365
$ws.Range("S6").Value = 0.003302991343996701
$ws.Range("T6").Value = 0.003302991343996702

$ws.Range("M7").Value = 3.867218333333334
$ws.Range("N7").Value = 11.601655
$ws.Range("O7").Value = 0.1566152977872902
$ws.Range("P7").Value = 0.1566152977872902
$ws.Range("Q7").Value = 22.73596439885334
$ws.Range("R7").Value = 204.62367958968
$ws.Range("S7").Value = 0.1524684818732999
$ws.Range("T7").Value = 0.1524684818732999

$ws.Range("N8").Value = 33.813685
$ws.Range("O8").Value = 0.4564642152831324
$ws.Range("P8").Value = 0.4564642152831324
$ws.Range("Q8").Value = 66.26526459837334
$ws.Range("R8").Value = 596.38738138536
$ws.Range("S8").Value = 0.4443780838588954
$ws.Range("T8").Value = 0.4443780838588955

$ws.Range("M9").Value = 5.654344666666667
$ws.Range("N9").Value = 16.963034
$ws.Range("O9").Value = 0.2289906587711778
$ws.Range("P9").Value = 0.2289906587711778
$ws.Range("Q9").Value = 33.24275175572267
$ws.Range("R9").Value = 299.184765801504
$ws.Range("S9").Value = 0.2229275083550726
$ws.Range("T9").Value = 0.2229275083550727

$ws.Range("M10").Value = 0.819389
$ws.Range("N10").Value = 2.458167
$ws.Range("O10").Value = 0.03318376186120772
$ws.Range("P10").Value = 0.03318376186120772
$ws.Range("Q10").Value = 4.817312478128001
$ws.Range("R10").Value = 43.355812303152
$ws.Range("S10").Value = 0.03230513152485952
$ws.Range("T10").Value = 0.03230513152485952

$ws.Range("M11").Value = 3.080288333333333
$ws.Range("N11").Value = 9.240865
$ws.Range("O11").Value = 0.1247460662971919
$ws.Range("P11").Value = 0.1247460662971919
$ws.Range("Q11").Value = 18.10948331549334
$ws.Range("R11").Value = 162.98534983944
$ws.Range("S11").Value = 0.1214430749531952
$ws.Range("T11").Value = 0.1214430749531952
